$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-05-12 06:00"

$ws.Range("H2").Value = $newTimestamp
$ws.Range("H3").Value = $newTimestamp
$ws.Range("H4").Value = $newTimestamp
$ws.Range("H5").Value = $newTimestamp
